$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("btmigrate_work")
$ws1.Range("D4").Value = "10.58.59.104"
$ws1.Range("L6").Value = "thynet.thy.com"
$ws1.Range("D7").Value = "Server5"
